$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K header value (year 2020), copying style of J4 (s=6)
$ws.Range("K4").Value = 2020
$ws.Range("K4").Style = $ws.Range("J4").Style

# New column K data values, copying values from column J, with a new style
# that replicates J's style but with applyNumberFormat explicitly set.
$ws.Range("K3").Value = $null
$ws.Range("K5").Value = 0.86
$ws.Range("K6").Value = 1.07
$ws.Range("K7").Value = 25.27
$ws.Range("K8").Value = 14
$ws.Range("K9").Value = 0.12
$ws.Range("K10").Value = 21.74
$ws.Range("K11").Value = 9.4600000000000009

# Selection to match the diff
$ws.Range("P7").Select()
